# Add a new "buscar" worksheet after "cadastro" with a small lookup table,
# make it the active sheet, and apply the light formatting used in the
# original edit (bold-ish header look via a dedicated style, grey body text,
# vertically centred, column auto width).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet, placed right after "cadastro"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "buscar"

# Header row
$ws2.Range("A1").Value = "Modelo"
$ws2.Range("B1").Value = "Resultado"

# Data rows
$ws2.Range("A2").Value = "HP CHROMEBOOK 14 G1 (ENERGY STAR)"
$ws2.Range("A3").Value = "HP CHROMEBOOK 14 G1 (ES)"
$ws2.Range("A4").Value = "HP ENVY - PORTÁTIL TOQUE 17T"

# Formatting: grey font + vertical centering for the model list
$ws2.Range("A2:A4").Font.Color = 2763306
$ws2.Range("A2:A4").VerticalAlignment = -4108

# A1 gets its own (otherwise plain) style
$ws2.Range("A1").WrapText = $false

# Column sizing
$ws2.Columns.Item(1).ColumnWidth = 36.28515625
$ws2.Columns.Item(2).ColumnWidth = 9.85546875

# Row 2 keeps a slightly smaller custom height, as in the source edit
$ws2.Rows.Item(2).RowHeight = 12.75

# Page setup / margins (metric defaults used by the original sheet)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1
$ws2.PageSetup.LeftMargin = 0.511811024 * 72
$ws2.PageSetup.RightMargin = 0.511811024 * 72
$ws2.PageSetup.TopMargin = 0.78740157499999996 * 72
$ws2.PageSetup.BottomMargin = 0.78740157499999996 * 72
$ws2.PageSetup.HeaderMargin = 0.31496062000000002 * 72
$ws2.PageSetup.FooterMargin = 0.31496062000000002 * 72

# Selection ends up one row below the data, matching the source edit
$ws2.Range("A5").Select() | Out-Null

# "buscar" is the active/visible sheet
$ws2.Activate()
